$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3  # E2: 2 -> 3
$ws.Cells.Item(2, 6).Value = 1  # F2: 0.6666666666666666 -> 1
$ws.Cells.Item(2, 7).Value = 21.08397466666667  # G2: 8.970950666666665 -> 21.08397466666667
$ws.Cells.Item(2, 8).Value = 63.251924  # H2: 26.912852 -> 63.251924
$ws.Cells.Item(2, 9).Value = 0.06331801375981215  # I2: 0.02838798528205506 -> 0.06331801375981215
$ws.Cells.Item(2, 10).Value = 0.06331801375981214  # J2: 0.02838798528205506 -> 0.06331801375981214
$ws.Cells.Item(2, 13).Value = 45.90594266666667  # M2: 31.22896466666667 -> 45.90594266666667
$ws.Cells.Item(2, 14).Value = 137.717828  # N2: 93.686894 -> 137.717828
$ws.Cells.Item(2, 15).Value = 0.3954672001633582  # O2: 0.2877106972998645 -> 0.3954672001633582
$ws.Cells.Item(2, 16).Value = 0.3954672001633583  # P2: 0.2877106972998645 -> 0.3954672001633583
$ws.Cells.Item(2, 17).Value = 967.8797322334525  # Q2: 280.153501395743 -> 967.8797322334525
$ws.Cells.Item(2, 18).Value = 8710.917590101071  # R2: 2521.381512561687 -> 8710.917590101071
$ws.Cells.Item(2, 19).Value = 0.0250401976214979  # S2: 0.008167527040438355 -> 0.0250401976214979
$ws.Cells.Item(2, 20).Value = 0.0250401976214979  # T2: 0.008167527040438355 -> 0.0250401976214979
$ws.Cells.Item(3, 5).Value = 3  # E3: 2 -> 3
$ws.Cells.Item(3, 6).Value = 1  # F3: 0.6666666666666666 -> 1
$ws.Cells.Item(3, 7).Value = 21.08397466666667  # G3: 8.970950666666665 -> 21.08397466666667
$ws.Cells.Item(3, 8).Value = 63.251924  # H3: 26.912852 -> 63.251924
$ws.Cells.Item(3, 9).Value = 0.06331801375981215  # I3: 0.02838798528205506 -> 0.06331801375981215
$ws.Cells.Item(3, 10).Value = 0.06331801375981214  # J3: 0.02838798528205506 -> 0.06331801375981214
$ws.Cells.Item(3, 15).Value = 0.3484294080560655  # O3: 0.3726247238124505 -> 0.3484294080560655
$ws.Cells.Item(3, 16).Value = 0.3484294080560656  # P3: 0.3726247238124506 -> 0.3484294080560656
$ws.Cells.Item(3, 17).Value = 852.7578571175058  # Q3: 362.8371209773885 -> 852.7578571175058
$ws.Cells.Item(3, 18).Value = 7674.820714057551  # R3: 3265.534088796496 -> 7674.820714057551
$ws.Cells.Item(3, 19).Value = 0.02206185805361716  # S3: 0.01057806517531768 -> 0.02206185805361716
$ws.Cells.Item(3, 20).Value = 0.02206185805361716  # T3: 0.01057806517531768 -> 0.02206185805361716
$ws.Cells.Item(4, 5).Value = 3  # E4: 2 -> 3
$ws.Cells.Item(4, 6).Value = 1  # F4: 0.6666666666666666 -> 1
$ws.Cells.Item(4, 7).Value = 21.08397466666667  # G4: 8.970950666666665 -> 21.08397466666667
$ws.Cells.Item(4, 8).Value = 63.251924  # H4: 26.912852 -> 63.251924
$ws.Cells.Item(4, 9).Value = 0.06331801375981215  # I4: 0.02838798528205506 -> 0.06331801375981215
$ws.Cells.Item(4, 10).Value = 0.06331801375981214  # J4: 0.02838798528205506 -> 0.06331801375981214
$ws.Cells.Item(4, 13).Value = 12.761795  # M4: 11.49855033333333 -> 12.761795
$ws.Cells.Item(4, 14).Value = 38.28538500000001  # N4: 34.495651 -> 38.28538500000001
$ws.Cells.Item(4, 15).Value = 0.1099393900775594  # O4: 0.1059354983315251 -> 0.1099393900775594
$ws.Cells.Item(4, 16).Value = 0.1099393900775594  # P4: 0.1059354983315251 -> 0.1099393900775594
$ws.Cells.Item(4, 17).Value = 269.0693624811934  # Q4: 103.1529277785169 -> 269.0693624811934
$ws.Cells.Item(4, 18).Value = 2421.62426233074  # R4: 928.3763500066519 -> 2421.62426233074
$ws.Cells.Item(4, 19).Value = 0.006961143813676261  # S4: 0.003007295367482504 -> 0.006961143813676261
$ws.Cells.Item(4, 20).Value = 0.006961143813676262  # T4: 0.003007295367482504 -> 0.006961143813676262
$ws.Cells.Item(5, 5).Value = 3  # E5: 2 -> 3
$ws.Cells.Item(5, 6).Value = 1  # F5: 0.6666666666666666 -> 1
$ws.Cells.Item(5, 7).Value = 21.08397466666667  # G5: 8.970950666666665 -> 21.08397466666667
$ws.Cells.Item(5, 8).Value = 63.251924  # H5: 26.912852 -> 63.251924
$ws.Cells.Item(5, 9).Value = 0.06331801375981215  # I5: 0.02838798528205506 -> 0.06331801375981215
$ws.Cells.Item(5, 10).Value = 0.06331801375981214  # J5: 0.02838798528205506 -> 0.06331801375981214
$ws.Cells.Item(5, 13).Value = 16.966758  # M5: 25.36964133333333 -> 16.966758
$ws.Cells.Item(5, 14).Value = 50.900274  # N5: 76.108924 -> 50.900274
$ws.Cells.Item(5, 15).Value = 0.1461640017030168  # O5: 0.2337290805561597 -> 0.1461640017030168
$ws.Cells.Item(5, 16).Value = 0.1461640017030168  # P5: 0.2337290805561597 -> 0.1461640017030168
$ws.Cells.Item(5, 17).Value = 357.726695847464  # Q5: 227.5898008323609 -> 357.726695847464
$ws.Cells.Item(5, 18).Value = 3219.540262627176  # R5: 2048.308207491248 -> 3219.540262627176
$ws.Cells.Item(5, 19).Value = 0.009254814271020824  # S5: 0.006635097698816524 -> 0.009254814271020824
$ws.Cells.Item(5, 20).Value = 0.009254814271020824  # T5: 0.006635097698816525 -> 0.009254814271020824
$ws.Cells.Item(6, 9).Value = 0.8174956765497907  # I6: 0.8614037742994388 -> 0.8174956765497907
$ws.Cells.Item(6, 10).Value = 0.8174956765497907  # J6: 0.8614037742994389 -> 0.8174956765497907
$ws.Cells.Item(6, 13).Value = 45.90594266666667  # M6: 31.22896466666667 -> 45.90594266666667
$ws.Cells.Item(6, 14).Value = 137.717828  # N6: 93.686894 -> 137.717828
$ws.Cells.Item(6, 15).Value = 0.3954672001633582  # O6: 0.2877106972998645 -> 0.3954672001633582
$ws.Cells.Item(6, 16).Value = 0.3954672001633583  # P6: 0.2877106972998645 -> 0.3954672001633583
$ws.Cells.Item(6, 17).Value = 12496.24632134645  # Q6: 8500.965499585682 -> 12496.24632134645
$ws.Cells.Item(6, 18).Value = 112466.2168921181  # R6: 76508.68949627114 -> 112466.2168921181
$ws.Cells.Item(6, 19).Value = 0.323292726350796  # S6: 0.2478350805604267 -> 0.323292726350796
$ws.Cells.Item(6, 20).Value = 0.3232927263507961  # T6: 0.2478350805604267 -> 0.3232927263507961
$ws.Cells.Item(7, 9).Value = 0.8174956765497907  # I7: 0.8614037742994388 -> 0.8174956765497907
$ws.Cells.Item(7, 10).Value = 0.8174956765497907  # J7: 0.8614037742994389 -> 0.8174956765497907
$ws.Cells.Item(7, 15).Value = 0.3484294080560655  # O7: 0.3726247238124505 -> 0.3484294080560655
$ws.Cells.Item(7, 16).Value = 0.3484294080560656  # P7: 0.3726247238124506 -> 0.3484294080560656
$ws.Cells.Item(7, 17).Value = 11009.91360818538  # Q7: 11009.91360818539 -> 11009.91360818538
$ws.Cells.Item(7, 18).Value = 99089.22247366846  # R7: 99089.22247366849 -> 99089.22247366846
$ws.Cells.Item(7, 19).Value = 0.2848395346686364  # S7: 0.3209803434893309 -> 0.2848395346686364
$ws.Cells.Item(7, 20).Value = 0.2848395346686364  # T7: 0.320980343489331 -> 0.2848395346686364
$ws.Cells.Item(8, 9).Value = 0.8174956765497907  # I8: 0.8614037742994388 -> 0.8174956765497907
$ws.Cells.Item(8, 10).Value = 0.8174956765497907  # J8: 0.8614037742994389 -> 0.8174956765497907
$ws.Cells.Item(8, 13).Value = 12.761795  # M8: 11.49855033333333 -> 12.761795
$ws.Cells.Item(8, 14).Value = 38.28538500000001  # N8: 34.495651 -> 38.28538500000001
$ws.Cells.Item(8, 15).Value = 0.1099393900775594  # O8: 0.1059354983315251 -> 0.1099393900775594
$ws.Cells.Item(8, 16).Value = 0.1099393900775594  # P8: 0.1059354983315251 -> 0.1099393900775594
$ws.Cells.Item(8, 17).Value = 3473.940944433009  # Q8: 3130.067894413794 -> 3473.940944433009
$ws.Cells.Item(8, 18).Value = 31265.46849989708  # R8: 28170.61104972415 -> 31265.46849989708
$ws.Cells.Item(8, 19).Value = 0.08987497607092577  # S8: 0.09125323809506763 -> 0.08987497607092577
$ws.Cells.Item(8, 20).Value = 0.08987497607092579  # T8: 0.09125323809506765 -> 0.08987497607092579
$ws.Cells.Item(9, 9).Value = 0.8174956765497907  # I9: 0.8614037742994388 -> 0.8174956765497907
$ws.Cells.Item(9, 10).Value = 0.8174956765497907  # J9: 0.8614037742994389 -> 0.8174956765497907
$ws.Cells.Item(9, 13).Value = 16.966758  # M9: 25.36964133333333 -> 16.966758
$ws.Cells.Item(9, 14).Value = 50.900274  # N9: 76.108924 -> 50.900274
$ws.Cells.Item(9, 15).Value = 0.1461640017030168  # O9: 0.2337290805561597 -> 0.1461640017030168
$ws.Cells.Item(9, 16).Value = 0.1461640017030168  # P9: 0.2337290805561597 -> 0.1461640017030168
$ws.Cells.Item(9, 17).Value = 4618.59129616847  # Q9: 6905.974886248109 -> 4618.59129616847
$ws.Cells.Item(9, 18).Value = 41567.32166551623  # R9: 62153.77397623299 -> 41567.32166551623
$ws.Cells.Item(9, 19).Value = 0.1194884394594325  # S9: 0.2013351121546136 -> 0.1194884394594325
$ws.Cells.Item(9, 20).Value = 0.1194884394594325  # T9: 0.2013351121546136 -> 0.1194884394594325
$ws.Cells.Item(10, 7).Value = 39.60693866666666  # G10: 34.42132633333333 -> 39.60693866666666
$ws.Cells.Item(10, 8).Value = 118.820816  # H10: 103.263979 -> 118.820816
$ws.Cells.Item(10, 9).Value = 0.1189449677837485  # I10: 0.1089240306459696 -> 0.1189449677837485
$ws.Cells.Item(10, 10).Value = 0.1189449677837485  # J10: 0.1089240306459696 -> 0.1189449677837485
$ws.Cells.Item(10, 13).Value = 45.90594266666667  # M10: 31.22896466666667 -> 45.90594266666667
$ws.Cells.Item(10, 14).Value = 137.717828  # N10: 93.686894 -> 137.717828
$ws.Cells.Item(10, 15).Value = 0.3954672001633582  # O10: 0.2877106972998645 -> 0.3954672001633582
$ws.Cells.Item(10, 16).Value = 0.3954672001633583  # P10: 0.2877106972998645 -> 0.3954672001633583
$ws.Cells.Item(10, 17).Value = 1818.193855634183  # Q10: 1074.94238384347 -> 1818.193855634183
$ws.Cells.Item(10, 18).Value = 16363.74470070765  # R10: 9674.481454591227 -> 16363.74470070765
$ws.Cells.Item(10, 19).Value = 0.04703883338295985  # S10: 0.03133860880986372 -> 0.04703883338295985
$ws.Cells.Item(10, 20).Value = 0.04703883338295985  # T10: 0.03133860880986372 -> 0.04703883338295985
$ws.Cells.Item(11, 7).Value = 39.60693866666666  # G11: 34.42132633333333 -> 39.60693866666666
$ws.Cells.Item(11, 8).Value = 118.820816  # H11: 103.263979 -> 118.820816
$ws.Cells.Item(11, 9).Value = 0.1189449677837485  # I11: 0.1089240306459696 -> 0.1189449677837485
$ws.Cells.Item(11, 10).Value = 0.1189449677837485  # J11: 0.1089240306459696 -> 0.1189449677837485
$ws.Cells.Item(11, 15).Value = 0.3484294080560655  # O11: 0.3726247238124505 -> 0.3484294080560655
$ws.Cells.Item(11, 16).Value = 0.3484294080560656  # P11: 0.3726247238124506 -> 0.3484294080560656
$ws.Cells.Item(11, 17).Value = 1601.933633403996  # Q11: 1392.19748397641 -> 1601.933633403996
$ws.Cells.Item(11, 18).Value = 14417.40270063597  # R11: 12529.7773557877 -> 14417.40270063597
$ws.Cells.Item(11, 19).Value = 0.04144392471613927  # S11: 0.04058778683599331 -> 0.04144392471613927
$ws.Cells.Item(11, 20).Value = 0.04144392471613927  # T11: 0.04058778683599332 -> 0.04144392471613927
$ws.Cells.Item(12, 7).Value = 39.60693866666666  # G12: 34.42132633333333 -> 39.60693866666666
$ws.Cells.Item(12, 8).Value = 118.820816  # H12: 103.263979 -> 118.820816
$ws.Cells.Item(12, 9).Value = 0.1189449677837485  # I12: 0.1089240306459696 -> 0.1189449677837485
$ws.Cells.Item(12, 10).Value = 0.1189449677837485  # J12: 0.1089240306459696 -> 0.1189449677837485
$ws.Cells.Item(12, 13).Value = 12.761795  # M12: 11.49855033333333 -> 12.761795
$ws.Cells.Item(12, 14).Value = 38.28538500000001  # N12: 34.495651 -> 38.28538500000001
$ws.Cells.Item(12, 15).Value = 0.1099393900775594  # O12: 0.1059354983315251 -> 0.1099393900775594
$ws.Cells.Item(12, 16).Value = 0.1099393900775594  # P12: 0.1059354983315251 -> 0.1099393900775594
$ws.Cells.Item(12, 17).Value = 505.4556318415733  # Q12: 395.7953533839255 -> 505.4556318415733
$ws.Cells.Item(12, 18).Value = 4549.10068657416  # R12: 3562.158180455329 -> 4549.10068657416
$ws.Cells.Item(12, 19).Value = 0.01307673721094026  # S12: 0.0115389214667591 -> 0.01307673721094026
$ws.Cells.Item(12, 20).Value = 0.01307673721094026  # T12: 0.0115389214667591 -> 0.01307673721094026
$ws.Cells.Item(13, 7).Value = 39.60693866666666  # G13: 34.42132633333333 -> 39.60693866666666
$ws.Cells.Item(13, 8).Value = 118.820816  # H13: 103.263979 -> 118.820816
$ws.Cells.Item(13, 9).Value = 0.1189449677837485  # I13: 0.1089240306459696 -> 0.1189449677837485
$ws.Cells.Item(13, 10).Value = 0.1189449677837485  # J13: 0.1089240306459696 -> 0.1189449677837485
$ws.Cells.Item(13, 13).Value = 16.966758  # M13: 25.36964133333333 -> 16.966758
$ws.Cells.Item(13, 14).Value = 50.900274  # N13: 76.108924 -> 50.900274
$ws.Cells.Item(13, 15).Value = 0.1461640017030168  # O13: 0.2337290805561597 -> 0.1461640017030168
$ws.Cells.Item(13, 16).Value = 0.1461640017030168  # P13: 0.2337290805561597 -> 0.1461640017030168
$ws.Cells.Item(13, 17).Value = 672.001343478176  # Q13: 873.2567032942885 -> 672.001343478176
$ws.Cells.Item(13, 18).Value = 6048.012091303583  # R13: 7859.310329648597 -> 6048.012091303583
$ws.Cells.Item(13, 19).Value = 0.01738547247370909  # S13: 0.02545871353335343 -> 0.01738547247370909
$ws.Cells.Item(13, 20).Value = 0.01738547247370909  # T13: 0.02545871353335344 -> 0.01738547247370909
$ws.Cells.Item(14, 7).Value = 0.08036333333333333  # G14: 0.405826 -> 0.08036333333333333
$ws.Cells.Item(14, 8).Value = 0.24109  # H14: 1.217478 -> 0.24109
$ws.Cells.Item(14, 9).Value = 0.0002413419066486121  # I14: 0.001284209772536402 -> 0.0002413419066486121
$ws.Cells.Item(14, 10).Value = 0.0002413419066486121  # J14: 0.001284209772536401 -> 0.0002413419066486121
$ws.Cells.Item(14, 13).Value = 45.90594266666667  # M14: 31.22896466666667 -> 45.90594266666667
$ws.Cells.Item(14, 14).Value = 137.717828  # N14: 93.686894 -> 137.717828
$ws.Cells.Item(14, 15).Value = 0.3954672001633582  # O14: 0.2877106972998645 -> 0.3954672001633582
$ws.Cells.Item(14, 16).Value = 0.3954672001633583  # P14: 0.2877106972998645 -> 0.3954672001633583
$ws.Cells.Item(14, 17).Value = 3.689154572502222  # Q14: 12.67352581481467 -> 3.689154572502222
$ws.Cells.Item(14, 18).Value = 33.20239115252  # R14: 114.061732333332 -> 33.20239115252
$ws.Cells.Item(14, 19).Value = [double]"9.544280810441319E-05"  # S14: 0.0003694808891357485 -> 9.544280810441319E-05
$ws.Cells.Item(14, 20).Value = [double]"9.544280810441322E-05"  # T14: 0.0003694808891357484 -> 9.544280810441322E-05
$ws.Cells.Item(15, 7).Value = 0.08036333333333333  # G15: 0.405826 -> 0.08036333333333333
$ws.Cells.Item(15, 8).Value = 0.24109  # H15: 1.217478 -> 0.24109
$ws.Cells.Item(15, 9).Value = 0.0002413419066486121  # I15: 0.001284209772536402 -> 0.0002413419066486121
$ws.Cells.Item(15, 10).Value = 0.0002413419066486121  # J15: 0.001284209772536401 -> 0.0002413419066486121
$ws.Cells.Item(15, 15).Value = 0.3484294080560655  # O15: 0.3726247238124505 -> 0.3484294080560655
$ws.Cells.Item(15, 16).Value = 0.3484294080560656  # P15: 0.3726247238124506 -> 0.3484294080560656
$ws.Cells.Item(15, 17).Value = 3.250357914368889  # Q15: 16.41395019648267 -> 3.250357914368889
$ws.Cells.Item(15, 18).Value = 29.25322122932  # R15: 147.725551768344 -> 29.25322122932
$ws.Cells.Item(15, 19).Value = [double]"8.409061767269815E-05"  # S15: 0.0004785283118086265 -> 8.409061767269815E-05
$ws.Cells.Item(15, 20).Value = [double]"8.409061767269816E-05"  # T15: 0.0004785283118086265 -> 8.409061767269816E-05
$ws.Cells.Item(16, 7).Value = 0.08036333333333333  # G16: 0.405826 -> 0.08036333333333333
$ws.Cells.Item(16, 8).Value = 0.24109  # H16: 1.217478 -> 0.24109
$ws.Cells.Item(16, 9).Value = 0.0002413419066486121  # I16: 0.001284209772536402 -> 0.0002413419066486121
$ws.Cells.Item(16, 10).Value = 0.0002413419066486121  # J16: 0.001284209772536401 -> 0.0002413419066486121
$ws.Cells.Item(16, 13).Value = 12.761795  # M16: 11.49855033333333 -> 12.761795
$ws.Cells.Item(16, 14).Value = 38.28538500000001  # N16: 34.495651 -> 38.28538500000001
$ws.Cells.Item(16, 15).Value = 0.1099393900775594  # O16: 0.1059354983315251 -> 0.1099393900775594
$ws.Cells.Item(16, 16).Value = 0.1099393900775594  # P16: 0.1059354983315251 -> 0.1099393900775594
$ws.Cells.Item(16, 17).Value = 1.025580385516667  # Q16: 4.666410687575334 -> 1.025580385516667
$ws.Cells.Item(16, 18).Value = 9.230223469650001  # R16: 41.997696188178 -> 9.230223469650001
$ws.Cells.Item(16, 19).Value = [double]"2.653298201710369E-05"  # S16: 0.0001360434022158582 -> 2.653298201710369E-05
$ws.Cells.Item(16, 20).Value = [double]"2.65329820171037E-05"  # T16: 0.0001360434022158582 -> 2.65329820171037E-05
$ws.Cells.Item(17, 7).Value = 0.08036333333333333  # G17: 0.405826 -> 0.08036333333333333
$ws.Cells.Item(17, 8).Value = 0.24109  # H17: 1.217478 -> 0.24109
$ws.Cells.Item(17, 9).Value = 0.0002413419066486121  # I17: 0.001284209772536402 -> 0.0002413419066486121
$ws.Cells.Item(17, 10).Value = 0.0002413419066486121  # J17: 0.001284209772536401 -> 0.0002413419066486121
$ws.Cells.Item(17, 13).Value = 16.966758  # M17: 25.36964133333333 -> 16.966758
$ws.Cells.Item(17, 14).Value = 50.900274  # N17: 76.108924 -> 50.900274
$ws.Cells.Item(17, 15).Value = 0.1461640017030168  # O17: 0.2337290805561597 -> 0.1461640017030168
$ws.Cells.Item(17, 16).Value = 0.1461640017030168  # P17: 0.2337290805561597 -> 0.1461640017030168
$ws.Cells.Item(17, 17).Value = 1.36350522874  # Q17: 10.29566006374133 -> 1.36350522874
$ws.Cells.Item(17, 18).Value = 12.27154705866  # R17: 92.66094057367199 -> 12.27154705866
$ws.Cells.Item(17, 19).Value = [double]"3.527549885439706E-05"  # S17: 0.0003001571693761681 -> 3.527549885439706E-05
$ws.Cells.Item(17, 20).Value = [double]"3.527549885439707E-05"  # T17: 0.0003001571693761681 -> 3.527549885439707E-05
